$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (was sitting in the empty,
#    numbered placeholder paragraph right before "#CONTRACTOR").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the "#REMARKSLIST" paragraph and the empty (sz=12)
#    paragraph that immediately precedes it, then:
#      - delete the whole empty sz=12 paragraph
#      - drop the leading tab run in the #REMARKSLIST paragraph
#      - plant the "_GoBack" bookmark at the start of that paragraph
# ------------------------------------------------------------------
$remarksIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*#REMARKSLIST*") {
        $remarksIndex = $i
    }
}

# Remove the empty placeholder paragraph (font size 6pt / w:sz 12) that
# sits directly above the remarks paragraph.
$emptyPara = $d.Paragraphs.Item($remarksIndex - 1)
if ($emptyPara.Range.Font.Size -eq 6) {
    $emptyPara.Range.Delete()
    $remarksIndex = $remarksIndex - 1
}

$remarksPara = $d.Paragraphs.Item($remarksIndex)

# Strip the leading tab character (the "<w:r><w:tab/></w:r>" run) from
# the front of the remarks paragraph, if present.
$firstChar = $d.Range($remarksPara.Range.Start, $remarksPara.Range.Start + 1)
if ($firstChar.Text -eq [char]9) {
    $firstChar.Delete()
}

# Re-create "_GoBack" collapsed at the very start of the paragraph.
$remarksPara = $d.Paragraphs.Item($remarksIndex)
$bmRange = $remarksPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 3) Update the cached "PAGE" field result in the header from 2 -> 1.
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdr.Range.Find.Execute("2", $false, $false, $false, $false, $false, $true, 1, $false, "1", 1) | Out-Null
